$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-36 down to 20-37
$ws.Rows("19:19").Insert()

# Populate the new row 19 with the new data record
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 44533
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103006
$ws.Range("J19").Value = "Nectarín"
$ws.Range("K19").Value = "Artic Pride"
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 270
$ws.Range("N19").Value = 24000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 24500
$ws.Range("Q19").Value = "$/bandeja 18 kilos granel"
$ws.Range("R19").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S19").Value = 1361
$ws.Range("T19").Value = 18
